{"js": "// Replace the 100 arithmetic-problem strings in the single table (20 rows x 5\n// cols) with their updated values, in document order, while preserving each\n// cell's existing paragraph/run formatting (font, size, alignment, etc.).\nconst replacements = [[\"2+61=\", \"43+43=\"], [\"55-39=\", \"22+53=\"], [\"97-61=\", \"39-21=\"], [\"2+94=\", \"75-54=\"], [\"28+50=\", \"31-21=\"], [\"29+25=\", \"13+36=\"], [\"37+48=\", \"42-36=\"], [\"26+64=\", \"33+14=\"], [\"92-38=\", \"67-58=\"], [\"98-20=\", \"43+33=\"], [\"84-69=\", \"10+71=\"], [\"63-62=\", \"77-15=\"], [\"97-96=\", \"25+13=\"], [\"12+28=\", \"89-39=\"], [\"26+18=\", \"52+7=\"], [\"23+15=\", \"6+88=\"], [\"19+50=\", \"61-43=\"], [\"47-45=\", \"10+70=\"], [\"15-1=\", \"24+56=\"], [\"54-2=\", \"25-2=\"], [\"20+26=\", \"52-6=\"], [\"40+23=\", \"21-8=\"], [\"51-16=\", \"21+74=\"], [\"18+11=\", \"43+36=\"], [\"46-43=\", \"16-13=\"], [\"33+28=\", \"35+64=\"], [\"35-25=\", \"76-54=\"], [\"81+9=\", \"25+35=\"], [\"18-3=\", \"25+72=\"], [\"30+35=\", \"5+65=\"], [\"93+1=\", \"45+21=\"], [\"92-29=\", \"0+31=\"], [\"51+0=\", \"80-49=\"], [\"61-44=\", \"65-41=\"], [\"32+1=\", \"7+1=\"], [\"53-0=\", \"87-85=\"], [\"67-61=\", \"72-28=\"], [\"9+8=\", \"38+16=\"], [\"94-31=\", \"39+29=\"], [\"96-17=\", \"81-69=\"], [\"46+37=\", \"66+31=\"], [\"8-4=\", \"80-61=\"], [\"5+51=\", \"21+37=\"], [\"95-71=\", \"57-25=\"], [\"49-21=\", \"23+74=\"], [\"85+11=\", \"88-9=\"], [\"64+33=\", \"48+22=\"], [\"3+91=\", \"20+68=\"], [\"19+3=\", \"74-7=\"], [\"46+18=\", \"64-56=\"], [\"95-45=\", \"40-30=\"], [\"14+12=\", \"62-58=\"], [\"86-57=\", \"69-61=\"], [\"22+60=\", \"6+2=\"], [\"43+45=\", \"90-74=\"], [\"83+8=\", \"89+7=\"], [\"61+16=\", \"2+20=\"], [\"65+17=\", \"77+2=\"], [\"85-35=\", \"16-6=\"], [\"80+4=\", \"94-50=\"], [\"33+2=\", \"30+10=\"], [\"36-0=\", \"81-38=\"], [\"56-24=\", \"12+38=\"], [\"78-9=\", \"2+97=\"], [\"2+80=\", \"25-22=\"], [\"23+42=\", \"48+16=\"], [\"56-16=\", \"3+76=\"], [\"5+60=\", \"9-5=\"], [\"89-12=\", \"22+41=\"], [\"54+15=\", \"39+2=\"], [\"60-27=\", \"76+16=\"], [\"20+53=\", \"41-12=\"], [\"9+70=\", \"92+6=\"], [\"14+59=\", \"83-73=\"], [\"46+9=\", \"96-67=\"], [\"62+2=\", \"62-30=\"], [\"29+38=\", \"74+12=\"], [\"43-15=\", \"96-60=\"], [\"42+49=\", \"31-16=\"], [\"5+47=\", \"19+45=\"], [\"66+10=\", \"77+7=\"], [\"23+7=\", \"2+65=\"], [\"12+77=\", \"53-17=\"], [\"65-18=\", \"14+80=\"], [\"62-28=\", \"76-6=\"], [\"67-0=\", \"40+10=\"], [\"11+17=\", \"0+22=\"], [\"96-69=\", \"41+35=\"], [\"48-48=\", \"12+22=\"], [\"2+78=\", \"84-35=\"], [\"70-28=\", \"84-2=\"], [\"14+0=\", \"63-49=\"], [\"41+56=\", \"98-7=\"], [\"62-17=\", \"81-13=\"], [\"92-82=\", \"20+43=\"], [\"75+2=\", \"5+21=\"], [\"91-49=\", \"17+78=\"], [\"21+41=\", \"70+22=\"], [\"86-75=\", \"39+20=\"], [\"63-53=\", \"19+29=\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Flatten all cells in row-major (document) order and grab the single\n// paragraph range inside each cell so insertText can swap the text without\n// touching the surrounding run/paragraph formatting.\nconst cells = [];\nfor (const row of rows) {\n  for (const cell of row.cells.items) {\n    cells.push(cell);\n  }\n}\n\nif (cells.length !== replacements.length) {\n  throw new Error(\n    \"Expected \" + replacements.length + \" cells, found \" + cells.length\n  );\n}\n\nconst paragraphs = [];\nfor (const cell of cells) {\n  cell.body.paragraphs.load(\"items\");\n}\nawait context.sync();\nfor (const cell of cells) {\n  paragraphs.push(cell.body.paragraphs.items[0]);\n}\n\nfor (const p of paragraphs) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < cells.length; i++) {\n  const [before, after] = replacements[i];\n  const actual = paragraphs[i].text;\n  // Defensive: only swap text that still matches the expected \"before\"\n  // value; otherwise leave the cell untouched.\n  if (actual === before) {\n    const range = paragraphs[i].getRange();\n    range.insertText(after, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem strings in the single table (20 rows x 5\n# cols) with their updated values, in document order, while preserving each\n# cell's existing formatting (Cell.Range.Text assignment keeps the run/\n# paragraph properties already present; only the visible text changes).\n$replacements = @(\n    @(\"2+61=\", \"43+43=\"),\n    @(\"55-39=\", \"22+53=\"),\n    @(\"97-61=\", \"39-21=\"),\n    @(\"2+94=\", \"75-54=\"),\n    @(\"28+50=\", \"31-21=\"),\n    @(\"29+25=\", \"13+36=\"),\n    @(\"37+48=\", \"42-36=\"),\n    @(\"26+64=\", \"33+14=\"),\n    @(\"92-38=\", \"67-58=\"),\n    @(\"98-20=\", \"43+33=\"),\n    @(\"84-69=\", \"10+71=\"),\n    @(\"63-62=\", \"77-15=\"),\n    @(\"97-96=\", \"25+13=\"),\n    @(\"12+28=\", \"89-39=\"),\n    @(\"26+18=\", \"52+7=\"),\n    @(\"23+15=\", \"6+88=\"),\n    @(\"19+50=\", \"61-43=\"),\n    @(\"47-45=\", \"10+70=\"),\n    @(\"15-1=\", \"24+56=\"),\n    @(\"54-2=\", \"25-2=\"),\n    @(\"20+26=\", \"52-6=\"),\n    @(\"40+23=\", \"21-8=\"),\n    @(\"51-16=\", \"21+74=\"),\n    @(\"18+11=\", \"43+36=\"),\n    @(\"46-43=\", \"16-13=\"),\n    @(\"33+28=\", \"35+64=\"),\n    @(\"35-25=\", \"76-54=\"),\n    @(\"81+9=\", \"25+35=\"),\n    @(\"18-3=\", \"25+72=\"),\n    @(\"30+35=\", \"5+65=\"),\n    @(\"93+1=\", \"45+21=\"),\n    @(\"92-29=\", \"0+31=\"),\n    @(\"51+0=\", \"80-49=\"),\n    @(\"61-44=\", \"65-41=\"),\n    @(\"32+1=\", \"7+1=\"),\n    @(\"53-0=\", \"87-85=\"),\n    @(\"67-61=\", \"72-28=\"),\n    @(\"9+8=\", \"38+16=\"),\n    @(\"94-31=\", \"39+29=\"),\n    @(\"96-17=\", \"81-69=\"),\n    @(\"46+37=\", \"66+31=\"),\n    @(\"8-4=\", \"80-61=\"),\n    @(\"5+51=\", \"21+37=\"),\n    @(\"95-71=\", \"57-25=\"),\n    @(\"49-21=\", \"23+74=\"),\n    @(\"85+11=\", \"88-9=\"),\n    @(\"64+33=\", \"48+22=\"),\n    @(\"3+91=\", \"20+68=\"),\n    @(\"19+3=\", \"74-7=\"),\n    @(\"46+18=\", \"64-56=\"),\n    @(\"95-45=\", \"40-30=\"),\n    @(\"14+12=\", \"62-58=\"),\n    @(\"86-57=\", \"69-61=\"),\n    @(\"22+60=\", \"6+2=\"),\n    @(\"43+45=\", \"90-74=\"),\n    @(\"83+8=\", \"89+7=\"),\n    @(\"61+16=\", \"2+20=\"),\n    @(\"65+17=\", \"77+2=\"),\n    @(\"85-35=\", \"16-6=\"),\n    @(\"80+4=\", \"94-50=\"),\n    @(\"33+2=\", \"30+10=\"),\n    @(\"36-0=\", \"81-38=\"),\n    @(\"56-24=\", \"12+38=\"),\n    @(\"78-9=\", \"2+97=\"),\n    @(\"2+80=\", \"25-22=\"),\n    @(\"23+42=\", \"48+16=\"),\n    @(\"56-16=\", \"3+76=\"),\n    @(\"5+60=\", \"9-5=\"),\n    @(\"89-12=\", \"22+41=\"),\n    @(\"54+15=\", \"39+2=\"),\n    @(\"60-27=\", \"76+16=\"),\n    @(\"20+53=\", \"41-12=\"),\n    @(\"9+70=\", \"92+6=\"),\n    @(\"14+59=\", \"83-73=\"),\n    @(\"46+9=\", \"96-67=\"),\n    @(\"62+2=\", \"62-30=\"),\n    @(\"29+38=\", \"74+12=\"),\n    @(\"43-15=\", \"96-60=\"),\n    @(\"42+49=\", \"31-16=\"),\n    @(\"5+47=\", \"19+45=\"),\n    @(\"66+10=\", \"77+7=\"),\n    @(\"23+7=\", \"2+65=\"),\n    @(\"12+77=\", \"53-17=\"),\n    @(\"65-18=\", \"14+80=\"),\n    @(\"62-28=\", \"76-6=\"),\n    @(\"67-0=\", \"40+10=\"),\n    @(\"11+17=\", \"0+22=\"),\n    @(\"96-69=\", \"41+35=\"),\n    @(\"48-48=\", \"12+22=\"),\n    @(\"2+78=\", \"84-35=\"),\n    @(\"70-28=\", \"84-2=\"),\n    @(\"14+0=\", \"63-49=\"),\n    @(\"41+56=\", \"98-7=\"),\n    @(\"62-17=\", \"81-13=\"),\n    @(\"92-82=\", \"20+43=\"),\n    @(\"75+2=\", \"5+21=\"),\n    @(\"91-49=\", \"17+78=\"),\n    @(\"21+41=\", \"70+22=\"),\n    @(\"86-75=\", \"39+20=\"),\n    @(\"63-53=\", \"19+29=\")\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif (($rowCount * $colCount) -ne $replacements.Count) {\n    throw \"Expected $($replacements.Count) cells, found $($rowCount * $colCount)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $pair = $replacements[$i]\n        $before = $pair[0]\n        $after = $pair[1]\n\n        $cell = $table.Cell($r, $c)\n        $range = $cell.Range\n        $current = $range.Text.TrimEnd([char]13, [char]7)\n\n        # Defensive: only swap text that still matches the expected \"before\"\n        # value; otherwise leave the cell untouched.\n        if ($current -eq $before) {\n            $range.Text = $after\n        }\n\n        $i++\n    }\n}\n"}
